# Atualização dicionario para novo modelo
# The "horas_trabalhadas" attribute row (DECIMAL, "Quantidade de horas
# trabalhadas neste serviço") is removed from the data dictionary table;
# the rows below it shift up to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 holds the "horas_trabalhadas" attribute (A7:H7). Deleting the
# whole row shifts rows 8-10 up to 7-9 and updates the rest of the sheet
# (shared formatting rows, merged cells) accordingly.
$ws.Rows.Item(7).Delete() | Out-Null

# Restore the selection Excel leaves behind after this kind of edit.
$ws.Range("H21").Select() | Out-Null
